# Generate Report for Handback
# Updates the "Ready for handoff" status to "Handed back: in sync with en-US"
# across all three sheets, refreshes the handback timestamps for both
# locales, clears the now-stale "latest handback name" / "error detail"
# values that pointed at an out-of-date handback, and widens a few columns
# to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E & F need to widen to fit the longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.16796875
$overview.Columns.Item(6).ColumnWidth = 29.16796875

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Refresh the handback timestamps to the latest successful handback run.
$zhcn.Range("L2").Value = "2017-02-09 08:14:53"
$zhcn.Range("L3").Value = "2017-02-09 08:14:53"

# The handback is now in sync, so the stale handback-name / error-detail
# fields are cleared out.
$zhcn.Range("M2").Value = ""
$zhcn.Range("R2").Value = ""
$zhcn.Range("R3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.16796875
$zhcn.Columns.Item(13).ColumnWidth = 23.001302083333336
$zhcn.Columns.Item(18).ColumnWidth = 12.834635416666668

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("L2").Value = "2017-02-09 08:15:19"
$dede.Range("L3").Value = "2017-02-09 08:15:19"

$dede.Range("M2").Value = ""
$dede.Range("R2").Value = ""
$dede.Range("R3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.16796875
$dede.Columns.Item(13).ColumnWidth = 23.001302083333336
$dede.Columns.Item(18).ColumnWidth = 12.834635416666668

Write-Host "Handback report generated"
